$d = $word.ActiveDocument

# --- Step 1: Replace the 14 existing equation placeholders with solved equations ---
$replacements = @(
    @("__Mg + __Br₂ → 6MgBr₂", "6Na + 3F₂ → 6NaF"),
    @("__K + 3Cl₂ →__6KCl", "8Mg + 4O₂ → 8MgO"),
    @("2Al + __Br₂ →__2AlBr₃", "2K + Cl₂ → 2KCl"),
    @("__K + __F₂ → 6KF", "2Fe₂O₃ + 3C → 4Fe + 3CO₂"),
    @("8Li + __Br₂ → __LiBr", "12Fe + 9O₂ → 6Fe₂O₃"),
    @("8Fe + __O₂ → __Fe₂O₃", "6Al + 9Br₂ → 6AlBr₃"),
    @("__Li +__3F₂ → 6LiF", "2Sn + 4H₂SO₄ → 2SnSO₄ + 4H₂O + 2SO₂"),
    @("__Al + 6Cl₂ →__4AlCl₃", "8Li + 4Br₂ → 8LiBr"),
    @("__Na + 3Cl₂ → __NaCl", "6K + 3Br₂ → 6KBr"),
    @("__Na + F₂ → __NaF", "Sn + 2H₂SO₄ → SnSO₄ + 2H₂O + SO₂"),
    @("Sn + __H₂SO₄ → SnSO₄ +__2H₂O + SO₂", "4Li + 2Cl₂ → 4LiCl"),
    @("__Li + 3Cl₂ →__6LiCl", "8Na + 4Br₂ → 8NaBr"),
    @("__Li + I₂ →__2LiI", "2Na + I₂ → 2NaI"),
    @("__K + I₂ →__2KI", "8Li + 4Cl₂ → 8LiCl"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find/Replace failed for: $old"
    }
}

# --- Step 2: Append new equation paragraphs after the last existing paragraph ---
$newParagraphTexts = @(
    "CH₄ + 2O₂ → CO₂ + 2H₂O",
    "8Na + 4F₂ → 8NaF",
    "4Sn + 8H₂SO₄ → 4SnSO₄ + 8H₂O + 4SO₂",
    "Sn + 2H₂SO₄ → SnSO₄ + 2H₂O + SO₂",
    "2Fe₂O₃ + 3C → 4Fe + 3CO₂",
    "8Mg + 4Cl₂ → 8MgCl₂",
    "2Li + Cl₂ → 2LiCl",
    "8K + 4Cl₂ → 8KCl",
    "6Mg + 3Cl₂ → 6MgCl₂",
    "CH₄ + 2O₂ → CO₂ + 2H₂O",
    "4Na + 2Br₂ → 4NaBr",
    "4Na + 2Cl₂ → 4NaCl",
    "6H₂ + 3O₂ → 6H₂O",
    "4Na + 2F₂ → 4NaF",
    "4H₂ + 2O₂ → 4H₂O",
    "2AgNO₃ + 2LiOH → 2AgOH + 2LiNO₃",
    "2K + F₂ → 2KF",
    "6Na + 3Cl₂ → 6NaCl",
    "2Na + F₂ → 2NaF",
    "2CH₄ + 4O₂ → 2CO₂ + 4H₂O",
    "6Mg + 3Cl₂ → 6MgCl₂",
    "4Al + 6Cl₂ → 4AlCl₃",
    "4Al + 6Cl₂ → 4AlCl₃",
    "8Fe₂O₃ + 12C → 16Fe + 12CO₂",
    "6K + 3I₂ → 6KI",
    "6Na + 3Br₂ → 6NaBr",
    "6Na + 3Br₂ → 6NaBr",
    "2AgNO₃ + 2LiOH → 2AgOH + 2LiNO₃",
    "4Na + 2F₂ → 4NaF",
    "6Mg + 3Br₂ → 6MgBr₂",
    "4Fe₂O₃ + 6C → 8Fe + 6CO₂",
    "6Mg + 3Br₂ → 6MgBr₂",
    "6Li + 3Br₂ → 6LiBr",
    "4Fe + 3O₂ → 2Fe₂O₃",
    "4Li + 2Br₂ → 4LiBr",
    "12Fe + 9O₂ → 6Fe₂O₃",
)

foreach ($text in $newParagraphTexts) {
    $lastPara = $d.Paragraphs($d.Paragraphs.Count)
    $lastPara.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs($d.Paragraphs.Count)
    $newPara.Alignment = 1
    $start = $newPara.Range.Start
    $r = $d.Range($start, $start)
    $r.InsertAfter($text)
    $textRange = $d.Range($start, $start + $text.Length)
    $textRange.Font.Size = 14
    $endRange = $d.Range($start + $text.Length, $start + $text.Length)
    $endRange.InsertAfter([char]11)
}

Write-Output $d.Paragraphs.Count
